$d = $word.ActiveDocument

# The phrase that appears 16 times in the document (once per use-case's
# "Performance" non-functional requirement). Only the 7th, 8th and 9th
# occurrences (the ones belonging to the "available rides/events",
# "added to their selected rides/events" and "see all the rides/events
# they signed on for" use cases) are being reworded from "3 seconds" to
# "less than a second", split across several runs as if the text had been
# edited incrementally by hand. The needle is the *entire* existing run
# (including the leading ": ") so the matched range's boundaries line up
# exactly with that <w:r> - replacing it then can't bleed into the bold
# "Performance" run just before it.
$needle = ": The system should be sent and saved within 3 seconds."

# Runs to build for each of the three target paragraphs (text, $true if
# the run should be a separate <w:r>; xml:space is derived automatically
# from leading/trailing whitespace by the runtime).
$targets = @{
    7 = @(
        ": The system should be sent and saved within ",
        "less ",
        "than",
        " a",
        " second",
        "."
    )
    8 = @(
        ": The system should be sent and saved",
        " within",
        " ",
        "less than a",
        " second",
        "."
    )
    9 = @(
        ": The system should be sent and ",
        "saved within ",
        "less tha",
        "n a second",
        "."
    )
}

$rng = $d.Content
$rng.Start = 0
$occurrence = 0

while ($rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $occurrence = $occurrence + 1

    if ($targets.ContainsKey($occurrence)) {
        $runs = $targets[$occurrence]

        # Replace the matched range's text with the first run's text (this
        # keeps it as its own <w:r> instead of merging into the preceding
        # run, which is what happens if the range is emptied out first).
        $rng.Text = $runs[0]
        $rng.Collapse(0)

        # Insert each remaining run's text one-by-one, collapsing to the
        # end after each insert so every InsertAfter call produces its own
        # <w:r>.
        for ($i = 1; $i -lt $runs.Count; $i++) {
            $rng.InsertAfter($runs[$i])
            $rng.Collapse(0)
        }
    }

    # Advance past this match (or past the inserted replacement) before
    # searching again.
    $rng.Start = $rng.End
    $rng.End = $d.Content.End
}

Write-Output "Processed $occurrence occurrence(s) of the performance requirement text."
